$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New-generated "K" values (column G) computed from the regenerated
# simulation (replacing the old "Strike#" based values).
$kValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 2
    9 = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 1
    29 = 3
    30 = 1
    31 = 2
    32 = 1
    33 = 2
    34 = 2
    35 = 1
    36 = 0
    37 = 3
    38 = 0
    39 = 2
    40 = 0
    41 = 4
    42 = 2
    43 = 1
    44 = 0
    45 = 2
    46 = 1
    47 = 2
    48 = 0
    49 = 1
    50 = 0
    51 = 1
    52 = 3
    53 = 0
    54 = 1
    55 = 1
    56 = 0
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 2
    62 = 2
    63 = 3
    64 = 2
    65 = 2
    66 = 1
    67 = 1
    68 = 0
    69 = 0
    70 = 2
    71 = 0
    72 = 1
    73 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Count) cells in column G (K)"
